$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) stays formatted as text so numeric-looking
# strings like "1.00" or "66.714.90" are not coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '66.714.90'
$ws.Range('E2').Value = '  -1.91%  '
$ws.Range('D3').Value = '2.441.13'
$ws.Range('E3').Value = '  -3.42%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '576.72'
$ws.Range('E5').Value = '  -3.02%  '
$ws.Range('D6').Value = '165.05'
$ws.Range('E6').Value = '  -6.62%  '
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('D8').Value = '0.503'
$ws.Range('E8').Value = '  -5.08%  '
$ws.Range('D9').Value = '2.438.89'
$ws.Range('E9').Value = '  -3.55%  '
$ws.Range('D10').Value = '0.133'
$ws.Range('E10').Value = '  -5.72%  '
$ws.Range('D11').Value = '0.163'
$ws.Range('E11').Value = '  -0.96%  '
$ws.Range('D12').Value = '4.79'
$ws.Range('E12').Value = '  -6.18%  '
$ws.Range('D13').Value = '0.326'
$ws.Range('E13').Value = '  -5.60%  '
$ws.Range('D14').Value = '24.97'
$ws.Range('E14').Value = '  -6.62%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '2.857.11'
$ws.Range('E15').Value = '  -4.45%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '66.394.20'
$ws.Range('E16').Value = '  -2.26%  '
$ws.Range('D17').Value = '0.0000166'
$ws.Range('E17').Value = '  -7.32%  '
$ws.Range('D18').Value = '2.432.56'
$ws.Range('E18').Value = '  -4.74%  '
$ws.Range('D19').Value = '11.15'
$ws.Range('E19').Value = '  -4.30%  '
$ws.Range('D20').Value = '7.47'
$ws.Range('E20').Value = '  -7.15%  '
$ws.Range('D21').Value = '351.53'
$ws.Range('E21').Value = '  -4.25%  '
$ws.Range('D22').Value = '3.98'
$ws.Range('E22').Value = '  -5.20%  '
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').Value = '68.88'
$ws.Range('E24').Value = '  -2.92%  '
$ws.Range('D25').Value = '4.18'
$ws.Range('E25').Value = '  -10.87%  '
$ws.Range('D26').Value = '1.71'
$ws.Range('E26').Value = '  -11.50%  '
$ws.Range('D27').Value = '8.84'
$ws.Range('E27').Value = '  -12.82%  '
$ws.Range('E28').Value = '  -0.20%  '
$ws.Range('D29').Value = '2.542.76'
$ws.Range('E29').Value = '  -4.19%  '
$ws.Range('B30').Value = 'Bittensor'
$ws.Range('C30').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D30').Value = '503.88'
$ws.Range('E30').Value = '  -5.88%  '
$ws.Range('B31').Value = 'PEPE'
$ws.Range('C31').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D31').Value = '0.0₃0882'
$ws.Range('E31').Value = '  -11.43%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = '7.71'
$ws.Range('E32').Value = '  -7.39%  '
$ws.Range('D33').Value = '1.76'
$ws.Range('E33').Value = '  -6.62%  '
$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').Value = '1.20'
$ws.Range('E34').Value = '  -9.76%  '
$ws.Range('B35').Value = 'FirstDigitalUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('D36').Value = '157.58'
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('E37').Value = '  -13.11%  '
$ws.Range('D38').Value = '18.54'
$ws.Range('E38').Value = '  -0.79%  '
$ws.Range('D39').Value = '18.33'
$ws.Range('E39').Value = '  -2.52%  '
$ws.Range('D40').Value = '1.33'
$ws.Range('E40').Value = '  -9.29%  '
$ws.Range('D41').Value = '1.65'
$ws.Range('E41').Value = '  -7.89%  '
$ws.Range('D42').Value = '0.322'
$ws.Range('E42').Value = '  -9.12%  '
$ws.Range('D43').Value = '4.67'
$ws.Range('E43').Value = '  -10.10%  '
$ws.Range('D44').Value = '39.16'
$ws.Range('E44').Value = '  -1.81%  '
$ws.Range('D45').Value = '2.31'
$ws.Range('E45').Value = '  -9.05%  '
$ws.Range('D46').Value = '139.76'
$ws.Range('E46').Value = '  -5.11%  '
$ws.Range('D47').Value = '3.44'
$ws.Range('E47').Value = '  -7.84%  '
$ws.Range('D48').Value = '0.505'
$ws.Range('E48').Value = '  -9.64%  '
$ws.Range('D49').Value = '1.57'
$ws.Range('E49').Value = '  -8.23%  '
$ws.Range('D50').Value = '0.0722'
$ws.Range('E50').Value = '  -4.77%  '
$ws.Range('D51').Value = '0.578'
$ws.Range('E51').Value = '  -3.30%  '
